# Daily TGP (terminal gate pricing) update.
# Each section's most-recent-day block of rows shifts down to become the
# second (previous-day) block, and the top block is replaced with a new
# day's prices. The old 'oldest day' block that used to occupy the second
# slot is discarded. All target values below were computed directly from
# the authoritative diff, so cells are written directly (no cell reads,
# which avoids a COM interop quirk where reading `.Value` on a cell
# returns a property descriptor instead of the actual value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New South Wales (Sydney-Botany, Sydney-Silverwater, Newcastle) - new day
$ws.Cells.Item(8, 1).Value = 46058
$ws.Cells.Item(8, 4).Value = 158.72
$ws.Cells.Item(8, 5).Value = 148.93
$ws.Cells.Item(8, 6).Value = 158.93
$ws.Cells.Item(8, 7).Value = 148.82

$ws.Cells.Item(9, 1).Value = 46058
$ws.Cells.Item(9, 4).Value = 158.72
$ws.Cells.Item(9, 5).Value = 148.93
$ws.Cells.Item(9, 6).Value = 158.93
$ws.Cells.Item(9, 7).Value = 148.82

$ws.Cells.Item(10, 1).Value = 46058
$ws.Cells.Item(10, 4).Value = 160.56
$ws.Cells.Item(10, 5).Value = 150.67
$ws.Cells.Item(10, 6).Value = 160.67
$ws.Cells.Item(10, 7).Value = 150.91

# New South Wales - shifted from previous day
$ws.Cells.Item(11, 1).Value = 46057
$ws.Cells.Item(11, 4).Value = 158.64
$ws.Cells.Item(11, 5).Value = 149.27
$ws.Cells.Item(11, 6).Value = 159.27
$ws.Cells.Item(11, 7).Value = 149.16

$ws.Cells.Item(12, 1).Value = 46057
$ws.Cells.Item(12, 4).Value = 158.64
$ws.Cells.Item(12, 5).Value = 149.27
$ws.Cells.Item(12, 6).Value = 159.27
$ws.Cells.Item(12, 7).Value = 149.16

$ws.Cells.Item(13, 1).Value = 46057
$ws.Cells.Item(13, 4).Value = 160.33
$ws.Cells.Item(13, 5).Value = 150.74
$ws.Cells.Item(13, 6).Value = 160.74
$ws.Cells.Item(13, 7).Value = 150.99

# Northern Territory (Darwin) - new day
$ws.Cells.Item(17, 1).Value = 46058
$ws.Cells.Item(17, 4).Value = 164.18
$ws.Cells.Item(17, 5).Value = 154.04
$ws.Cells.Item(17, 6).Value = 164.04

# Northern Territory - shifted from previous day
$ws.Cells.Item(18, 1).Value = 46057
$ws.Cells.Item(18, 4).Value = 163.94
$ws.Cells.Item(18, 5).Value = 154.2
$ws.Cells.Item(18, 6).Value = 164.2

# Queensland (Brisbane, Cairns, Gladstone, Mackay, Townsville) - new day
$ws.Cells.Item(22, 1).Value = 46058
$ws.Cells.Item(22, 4).Value = 160.12
$ws.Cells.Item(22, 5).Value = 150.72
$ws.Cells.Item(22, 6).Value = 160.32
$ws.Cells.Item(22, 7).Value = 152.48

$ws.Cells.Item(23, 1).Value = 46058
$ws.Cells.Item(23, 4).Value = 165.33
$ws.Cells.Item(23, 5).Value = 156.68
$ws.Cells.Item(23, 6).Value = 166.68

$ws.Cells.Item(24, 1).Value = 46058
$ws.Cells.Item(24, 4).Value = 165.49
$ws.Cells.Item(24, 5).Value = 157.3
$ws.Cells.Item(24, 6).Value = 167.3

$ws.Cells.Item(25, 1).Value = 46058
$ws.Cells.Item(25, 4).Value = 165.49
$ws.Cells.Item(25, 5).Value = 156.83
$ws.Cells.Item(25, 6).Value = 166.83
$ws.Cells.Item(25, 7).Value = 157.69

$ws.Cells.Item(26, 1).Value = 46058
$ws.Cells.Item(26, 4).Value = 165.08
$ws.Cells.Item(26, 5).Value = 158.42
$ws.Cells.Item(26, 6).Value = 168.42

# Queensland - shifted from previous day
$ws.Cells.Item(27, 1).Value = 46057
$ws.Cells.Item(27, 4).Value = 159.82
$ws.Cells.Item(27, 5).Value = 150.73
$ws.Cells.Item(27, 6).Value = 160.33
$ws.Cells.Item(27, 7).Value = 152.49

$ws.Cells.Item(28, 1).Value = 46057
$ws.Cells.Item(28, 4).Value = 165.11
$ws.Cells.Item(28, 5).Value = 156.86
$ws.Cells.Item(28, 6).Value = 166.86

$ws.Cells.Item(29, 1).Value = 46057
$ws.Cells.Item(29, 4).Value = 165.27
$ws.Cells.Item(29, 5).Value = 157.5
$ws.Cells.Item(29, 6).Value = 167.5

$ws.Cells.Item(30, 1).Value = 46057
$ws.Cells.Item(30, 4).Value = 165.27
$ws.Cells.Item(30, 5).Value = 157.03
$ws.Cells.Item(30, 6).Value = 167.03
$ws.Cells.Item(30, 7).Value = 157.88

$ws.Cells.Item(31, 1).Value = 46057
$ws.Cells.Item(31, 4).Value = 164.86
$ws.Cells.Item(31, 5).Value = 158.61
$ws.Cells.Item(31, 6).Value = 168.61

# South Australia (Adelaide) - new day
$ws.Cells.Item(35, 1).Value = 46058
$ws.Cells.Item(35, 4).Value = 158.5
$ws.Cells.Item(35, 5).Value = 148.36
$ws.Cells.Item(35, 6).Value = 157.36

# South Australia - shifted from previous day
$ws.Cells.Item(36, 1).Value = 46057
$ws.Cells.Item(36, 4).Value = 158.83
$ws.Cells.Item(36, 5).Value = 148.54
$ws.Cells.Item(36, 6).Value = 157.54

# Tasmania (Burnie, Hobart) - new day
$ws.Cells.Item(40, 1).Value = 46058
$ws.Cells.Item(40, 4).Value = 164.85
$ws.Cells.Item(40, 5).Value = 155.93
$ws.Cells.Item(40, 6).Value = 165.93

$ws.Cells.Item(41, 1).Value = 46058
$ws.Cells.Item(41, 4).Value = 164.57
$ws.Cells.Item(41, 5).Value = 156.35
$ws.Cells.Item(41, 6).Value = 166.35

# Tasmania - shifted from previous day
$ws.Cells.Item(42, 1).Value = 46057
$ws.Cells.Item(42, 4).Value = 164.65
$ws.Cells.Item(42, 5).Value = 156.04
$ws.Cells.Item(42, 6).Value = 166.04

$ws.Cells.Item(43, 1).Value = 46057
$ws.Cells.Item(43, 4).Value = 164.36
$ws.Cells.Item(43, 5).Value = 156.46
$ws.Cells.Item(43, 6).Value = 166.46

# Victoria (Geelong, Melbourne) - new day
$ws.Cells.Item(47, 1).Value = 46058
$ws.Cells.Item(47, 4).Value = 159.44
$ws.Cells.Item(47, 5).Value = 149.76
$ws.Cells.Item(47, 6).Value = 159.76

$ws.Cells.Item(48, 1).Value = 46058
$ws.Cells.Item(48, 4).Value = 159.08
$ws.Cells.Item(48, 5).Value = 149.71
$ws.Cells.Item(48, 6).Value = 159.71

# Victoria - shifted from previous day
$ws.Cells.Item(49, 1).Value = 46057
$ws.Cells.Item(49, 4).Value = 158.93
$ws.Cells.Item(49, 5).Value = 150.22
$ws.Cells.Item(49, 6).Value = 160.22

$ws.Cells.Item(50, 1).Value = 46057
$ws.Cells.Item(50, 4).Value = 158.55
$ws.Cells.Item(50, 5).Value = 150.16
$ws.Cells.Item(50, 6).Value = 160.16

# Western Australia (Broome, Esperance, Geraldton, Kalgoorlie, Perth, Port Hedland) - new day
$ws.Cells.Item(54, 1).Value = 46058
$ws.Cells.Item(54, 4).Value = 173.87
$ws.Cells.Item(54, 5).Value = 164.02
$ws.Cells.Item(54, 6).Value = 174.02

$ws.Cells.Item(55, 1).Value = 46058
$ws.Cells.Item(55, 4).Value = 163.29
$ws.Cells.Item(55, 5).Value = 162.06
$ws.Cells.Item(55, 6).Value = 172.06

$ws.Cells.Item(56, 1).Value = 46058
$ws.Cells.Item(56, 4).Value = 163.17

$ws.Cells.Item(57, 1).Value = 46058
$ws.Cells.Item(57, 4).Value = 163.78
$ws.Cells.Item(57, 5).Value = 156.49

$ws.Cells.Item(58, 1).Value = 46058
$ws.Cells.Item(58, 4).Value = 159.55
$ws.Cells.Item(58, 5).Value = 152.38
$ws.Cells.Item(58, 6).Value = 162.38

$ws.Cells.Item(59, 1).Value = 46058
$ws.Cells.Item(59, 4).Value = 166.57
$ws.Cells.Item(59, 5).Value = 162.24

# Western Australia - shifted from previous day
$ws.Cells.Item(60, 1).Value = 46057
$ws.Cells.Item(60, 4).Value = 173.65
$ws.Cells.Item(60, 5).Value = 164.28
$ws.Cells.Item(60, 6).Value = 174.28

$ws.Cells.Item(61, 1).Value = 46057
$ws.Cells.Item(61, 4).Value = 163.07
$ws.Cells.Item(61, 5).Value = 162.12
$ws.Cells.Item(61, 6).Value = 172.12

$ws.Cells.Item(62, 1).Value = 46057
$ws.Cells.Item(62, 4).Value = 162.95

$ws.Cells.Item(63, 1).Value = 46057
$ws.Cells.Item(63, 4).Value = 163.53
$ws.Cells.Item(63, 5).Value = 156.54

$ws.Cells.Item(64, 1).Value = 46057
$ws.Cells.Item(64, 4).Value = 159.3
$ws.Cells.Item(64, 5).Value = 152.44
$ws.Cells.Item(64, 6).Value = 162.44

$ws.Cells.Item(65, 1).Value = 46057
$ws.Cells.Item(65, 4).Value = 166.32
$ws.Cells.Item(65, 5).Value = 162.48
